# Auto-generated update of cryptos.xlsx data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2
Set-TextCell "D2" "46.601.15"
Set-TextCell "E2" "  +3.38%  "

# Row 3
Set-TextCell "D3" "2.274.09"
Set-TextCell "E3" "  +0.24%  "

# Row 4
Set-TextCell "D4" "0.999"
Set-TextCell "E4" "  -0.04%  "

# Row 5
Set-TextCell "D5" "301.36"
Set-TextCell "E5" "  -0.15%  "

# Row 6
Set-TextCell "D6" "99.97"
Set-TextCell "E6" "  +5.47%  "

# Row 7
Set-TextCell "D7" "0.562"
Set-TextCell "E7" "  -0.61%  "

# Row 8
Set-TextCell "D8" "1.00"
Set-TextCell "E8" "  +0.16%  "

# Row 9
Set-TextCell "D9" "0.518"
Set-TextCell "E9" "  +1.45%  "

# Row 10
Set-TextCell "D10" "35.99"
Set-TextCell "E10" "  +4.56%  "

# Row 11
Set-TextCell "D11" "0.0782"
Set-TextCell "E11" "  -1.03%  "

# Row 12
Set-TextCell "D12" "7.21"
Set-TextCell "E12" "  -0.10%  "

# Row 13
Set-TextCell "E13" "  -1.13%  "

# Row 14
Set-TextCell "D14" "2.622.28"
Set-TextCell "E14" "  +0.29%  "

# Row 15
Set-TextCell "D15" "2.275.29"
Set-TextCell "E15" "  +0.24%  "

# Row 16
Set-TextCell "D16" "13.67"
Set-TextCell "E16" "  -0.96%  "

# Row 17
Set-TextCell "D17" "0.802"
Set-TextCell "E17" "  +0.51%  "

# Row 18
Set-TextCell "D18" "46.569.99"
Set-TextCell "E18" "  +3.75%  "

# Row 19
Set-TextCell "D19" "13.08"
Set-TextCell "E19" "  +1.17%  "

# Row 20
Set-TextCell "D20" "0.0₃0930"
Set-TextCell "E20" "  +0.66%  "

# Row 21
Set-TextCell "D21" "5.93"
Set-TextCell "E21" "  -2.79%  "

# Row 22
Set-TextCell "D22" "65.31"
Set-TextCell "E22" "  -0.01%  "

# Row 23
Set-TextCell "D23" "247.03"
Set-TextCell "E23" "  +3.24%  "

# Row 24
Set-TextCell "D24" "2.86"
Set-TextCell "E24" "  -1.31%  "

# Row 25
Set-TextCell "E25" "  -0.02%  "

# Row 26
Set-TextCell "D26" "1.89"
Set-TextCell "E26" "  -1.55%  "

# Row 27
Set-TextCell "D27" "42.46"
Set-TextCell "E27" "  +3.40%  "

# Row 28
Set-TextCell "E28" "  -0.77%  "

# Row 29
Set-TextCell "D29" "9.79"
Set-TextCell "E29" "  +2.36%  "

# Row 30
Set-TextCell "D30" "19.85"
Set-TextCell "E30" "  +1.27%  "

# Row 31
Set-TextCell "D31" "2.78"
Set-TextCell "E31" "  +8.56%  "

# Row 32
Set-TextCell "D32" "5.50"
Set-TextCell "E32" "  -2.07%  "

# Row 33
Set-TextCell "D33" "146.66"
Set-TextCell "E33" "  -3.89%  "

# Row 34
Set-TextCell "B34" "LidoDAOToken"
Set-TextCell "C34" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextCell "D34" "3.27"
Set-TextCell "E34" "  +11.91%  "

# Row 35
Set-TextCell "B35" "Hedera"
Set-TextCell "C35" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D35" "0.0777"
Set-TextCell "E35" "  -1.34%  "

# Row 36
Set-TextCell "E36" "  +9.28%  "

# Row 37
Set-TextCell "D37" "0.116"
Set-TextCell "E37" "  -1.09%  "

# Row 38
Set-TextCell "D38" "15.88"
Set-TextCell "E38" "  +14.24%  "

# Row 39
Set-TextCell "D39" "1.74"
Set-TextCell "E39" "  +0.21%  "

# Row 40
Set-TextCell "D40" "3.96"
Set-TextCell "E40" "  +4.98%  "

# Row 41
Set-TextCell "D41" "3.27"
Set-TextCell "E41" "  +0.13%  "

# Row 42
Set-TextCell "D42" "0.0298"
Set-TextCell "E42" "  -2.32%  "

# Row 43
Set-TextCell "D43" "0.998"
Set-TextCell "E43" "  -0.13%  "

# Row 44
Set-TextCell "D44" "1.98"
Set-TextCell "E44" "  +2.48%  "

# Row 45
Set-TextCell "D45" "1.808.54"
Set-TextCell "E45" "  +1.28%  "

# Row 46
Set-TextCell "D46" "90.04"
Set-TextCell "E46" "  +18.76%  "

# Row 47
Set-TextCell "D47" "0.191"
Set-TextCell "E47" "  -0.42%  "

# Row 48
Set-TextCell "D48" "71.80"
Set-TextCell "E48" "  +1.78%  "

# Row 49
Set-TextCell "D49" "4.86"
Set-TextCell "E49" "  +4.26%  "

# Row 50
Set-TextCell "D50" "94.49"
Set-TextCell "E50" "  -2.42%  "

# Row 51
Set-TextCell "D51" "2.498.42"
Set-TextCell "E51" "  +0.40%  "
